# Refresh the cryptos list snapshot (Price + Volume(1h) columns) for rows 2-51.
# Mirrors a GitHub Actions scheduled data refresh: only column D (Price) and
# column E (Volume(1h)) change; everything else (rank, coin name, link) is untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "45.268.37"
$ws.Cells.Item(2, 5).Value = "  +5.06%  "
$ws.Cells.Item(3, 4).Value = "2.451.60"
$ws.Cells.Item(3, 5).Value = "  +3.44%  "
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "0.999"
$ws.Cells.Item(4, 5).Value = "  -0.08%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "319.23"
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "104.67"
$ws.Cells.Item(6, 5).Value = "  +9.11%  "
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.517"
$ws.Cells.Item(7, 5).Value = "  +2.69%  "
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.998"
$ws.Cells.Item(8, 5).Value = "  -0.21%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.534"
$ws.Cells.Item(9, 5).Value = "  +10.38%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "36.03"
$ws.Cells.Item(10, 5).Value = "  +4.60%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0807"
$ws.Cells.Item(11, 5).Value = "  +2.21%  "
$ws.Cells.Item(12, 5).Value = "  -2.62%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "18.48"
$ws.Cells.Item(13, 5).Value = "  +0.83%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "7.04"
$ws.Cells.Item(14, 5).Value = "  +3.34%  "
$ws.Cells.Item(15, 4).Value = "2.831.76"
$ws.Cells.Item(15, 5).Value = "  +3.54%  "
$ws.Cells.Item(16, 4).Value = "2.505.49"
$ws.Cells.Item(16, 5).Value = "  +6.23%  "
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "0.842"
$ws.Cells.Item(17, 5).Value = "  +4.69%  "
$ws.Cells.Item(18, 4).Value = "45.148.84"
$ws.Cells.Item(18, 5).Value = "  +4.74%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "12.39"
$ws.Cells.Item(19, 5).Value = "  +3.56%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "6.39"
$ws.Cells.Item(20, 5).Value = "  +1.46%  "
$ws.Cells.Item(21, 4).Value = "0.0₃0924"
$ws.Cells.Item(21, 5).Value = "  +3.84%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "69.22"
$ws.Cells.Item(22, 5).Value = "  +1.84%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "244.60"
$ws.Cells.Item(23, 5).Value = "  +3.88%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "2.30"
$ws.Cells.Item(24, 5).Value = "  +3.63%  "
$ws.Cells.Item(25, 5).Value = "  +3.13%  "
$ws.Cells.Item(26, 5).Value = "  +0.18%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "25.57"
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "2.19"
$ws.Cells.Item(28, 5).Value = "  -7.47%  "
$ws.Cells.Item(29, 5).Value = "  +2.61%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "34.10"
$ws.Cells.Item(30, 5).Value = "  +6.44%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "49.58"
$ws.Cells.Item(31, 5).Value = "  +3.40%  "
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "0.129"
$ws.Cells.Item(32, 5).Value = "  +15.75%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "20.17"
$ws.Cells.Item(33, 5).Value = "  +13.21%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "5.25"
$ws.Cells.Item(34, 5).Value = "  +3.90%  "
$ws.Cells.Item(35, 5).Value = "  +0.24%  "
$ws.Cells.Item(36, 5).Value = "  +3.81%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "1.92"
$ws.Cells.Item(37, 5).Value = "  +5.01%  "
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "4.51"
$ws.Cells.Item(38, 5).Value = "  +4.23%  "
$ws.Cells.Item(39, 5).Value = "  +1.04%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "125.15"
$ws.Cells.Item(40, 5).Value = "  -2.83%  "
$ws.Cells.Item(41, 5).Value = "  +2.44%  "
$ws.Cells.Item(42, 5).Value = "  -2.61%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "21.25"
$ws.Cells.Item(43, 5).Value = "  +0.28%  "
$ws.Cells.Item(44, 5).Value = "  +4.46%  "
$ws.Cells.Item(45, 4).Value = "1.948.78"
$ws.Cells.Item(45, 5).Value = "  +1.01%  "
$ws.Cells.Item(46, 5).Value = "  +7.67%  "
$ws.Cells.Item(47, 5).Value = "  -0.63%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "9.29"
$ws.Cells.Item(48, 5).Value = "  +0.89%  "
$ws.Cells.Item(49, 5).Value = "  +18.09%  "
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "76.32"
$ws.Cells.Item(50, 5).Value = "  +6.52%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "54.03"
$ws.Cells.Item(51, 5).Value = "  +4.52%  "
